# Auto-generated edit script: updates market-price derived columns (H:N)
# across multiple sheets to match the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 290.2
$ws.Range("I5").Value = 383.66666
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 383.66666
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -268.66666
$ws.Range("N5").Value = -380
$ws.Range("H18").Value = 222.8
$ws.Range("I18").Value = 222.8
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 222.8
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 61.19999999999999
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 595.9286
$ws.Range("I19").Value = 686.6
$ws.Range("J19").Value = 369.25
$ws.Range("K19").Value = 686.6
$ws.Range("L19").Value = 369.25
$ws.Range("M19").Value = -511.6
$ws.Range("N19").Value = -719.25
$ws.Range("H39").Value = 574.7273
$ws.Range("I39").Value = 110.625
$ws.Range("J39").Value = 839.9286
$ws.Range("K39").Value = 331.875
$ws.Range("L39").Value = 2519.7858
$ws.Range("M39").Value = -35.875
$ws.Range("N39").Value = -3111.7858
$ws.Range("H62").Value = 3718.4814
$ws.Range("I62").Value = 2661.111
$ws.Range("J62").Value = 5833.222
$ws.Range("K62").Value = 2661.111
$ws.Range("L62").Value = 5833.222
$ws.Range("M62").Value = -2037.111
$ws.Range("N62").Value = -7081.222
$ws.Range("H65").Value = 3718.4814
$ws.Range("I65").Value = 2661.111
$ws.Range("J65").Value = 5833.222
$ws.Range("K65").Value = 13305.555
$ws.Range("L65").Value = 29166.11
$ws.Range("M65").Value = -10185.555
$ws.Range("N65").Value = -35406.11
$ws.Range("H98").Value = 1214.9333
$ws.Range("I98").Value = 1254.2
$ws.Range("J98").Value = 1136.4
$ws.Range("K98").Value = 1254.2
$ws.Range("L98").Value = 1136.4
$ws.Range("M98").Value = 243.8
$ws.Range("N98").Value = -4132.4
$ws.Range("H112").Value = 1259.303
$ws.Range("J112").Value = 1298.4445
$ws.Range("L112").Value = 3895.3335
$ws.Range("N112").Value = -6111.333500000001
$ws.Range("H113").Value = 6484.6587
$ws.Range("I113").Value = 3248.4443
$ws.Range("J113").Value = 9017.348
$ws.Range("K113").Value = 3248.4443
$ws.Range("L113").Value = 9017.348
$ws.Range("M113").Value = 5.555699999999888
$ws.Range("N113").Value = -15525.348
$ws.Range("H122").Value = 1214.9333
$ws.Range("I122").Value = 1254.2
$ws.Range("J122").Value = 1136.4
$ws.Range("K122").Value = 3762.6
$ws.Range("L122").Value = 3409.2
$ws.Range("M122").Value = -1312.6
$ws.Range("N122").Value = -8309.200000000001
$ws.Range("H137").Value = 1716.8334
$ws.Range("I137").Value = 2183.6667
$ws.Range("J137").Value = 1250
$ws.Range("K137").Value = 6551.000100000001
$ws.Range("L137").Value = 3750
$ws.Range("M137").Value = -4001.000100000001
$ws.Range("N137").Value = -8850
$ws.Range("H138").Value = 3323.2
$ws.Range("I138").Value = 2160.7046
$ws.Range("J138").Value = 4236.5894
$ws.Range("K138").Value = 6482.1138
$ws.Range("L138").Value = 12709.7682
$ws.Range("M138").Value = -1342.1138
$ws.Range("N138").Value = -22989.7682

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 10319.138
$ws.Range("I32").Value = 7409.769
$ws.Range("J32").Value = 35533.668
$ws.Range("K32").Value = 7409.769
$ws.Range("L32").Value = 35533.668
$ws.Range("M32").Value = -7122.769
$ws.Range("N32").Value = -36107.668
$ws.Range("H33").Value = 10026
$ws.Range("I33").Value = 10026
$ws.Range("K33").Value = 10026
$ws.Range("M33").Value = -9697
$ws.Range("H45").Value = 2700.0667
$ws.Range("I45").Value = 2089.111
$ws.Range("J45").Value = 3616.5
$ws.Range("K45").Value = 2089.111
$ws.Range("L45").Value = 3616.5
$ws.Range("M45").Value = -1712.111
$ws.Range("N45").Value = -4370.5
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1942.75
$ws.Range("I74").Value = 1312.7273
$ws.Range("J74").Value = 2932.7856
$ws.Range("K74").Value = 1312.7273
$ws.Range("L74").Value = 2932.7856
$ws.Range("M74").Value = -438.7273
$ws.Range("N74").Value = -4680.7856
$ws.Range("H77").Value = 1942.75
$ws.Range("I77").Value = 1312.7273
$ws.Range("J77").Value = 2932.7856
$ws.Range("K77").Value = 6563.636500000001
$ws.Range("L77").Value = 14663.928
$ws.Range("M77").Value = -2195.636500000001
$ws.Range("N77").Value = -23399.928

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3830.8555
$ws.Range("I31").Value = 1955.7164
$ws.Range("J31").Value = 11683
$ws.Range("K31").Value = 1955.7164
$ws.Range("L31").Value = 11683
$ws.Range("M31").Value = -1660.7164
$ws.Range("N31").Value = -12273
$ws.Range("H34").Value = 3830.8555
$ws.Range("I34").Value = 1955.7164
$ws.Range("J34").Value = 11683
$ws.Range("K34").Value = 1955.7164
$ws.Range("L34").Value = 11683
$ws.Range("M34").Value = -1753.7164
$ws.Range("N34").Value = -12087
$ws.Range("H122").Value = 3394.6428
$ws.Range("I122").Value = 1201.1
$ws.Range("J122").Value = 8878.5
$ws.Range("K122").Value = 3603.3
$ws.Range("L122").Value = 26635.5
$ws.Range("M122").Value = -1153.3
$ws.Range("N122").Value = -31535.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2827
$ws.Range("N16").ClearContents()
$ws.Range("H103").Value = 2714.1667
$ws.Range("I103").Value = 71.25
$ws.Range("J103").Value = 8000
$ws.Range("K103").Value = 213.75
$ws.Range("L103").Value = 24000
$ws.Range("M103").Value = 665.25
$ws.Range("N103").Value = -25758
$ws.Range("H112").Value = 1640.1
$ws.Range("I112").Value = 967
$ws.Range("J112").Value = 1928.5714
$ws.Range("K112").Value = 2901
$ws.Range("L112").Value = 5785.7142
$ws.Range("M112").Value = -1793
$ws.Range("N112").Value = -8001.7142

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2331.25
$ws.Range("I97").Value = 2950
$ws.Range("K97").Value = 2950
$ws.Range("M97").Value = -2454
$ws.Range("H122").Value = 98556.16
$ws.Range("I122").Value = 148494.36
$ws.Range("J122").Value = 4228.4443
$ws.Range("K122").Value = 445483.08
$ws.Range("L122").Value = 12685.3329
$ws.Range("M122").Value = -443033.08
$ws.Range("N122").Value = -17585.3329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H40").Value = 57288.223
$ws.Range("I40").Value = 64074.25
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 64074.25
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -63938.25
$ws.Range("N40").Value = -3272
$ws.Range("H68").Value = 1695
$ws.Range("I68").Value = 1396
$ws.Range("J68").Value = 1894.3334
$ws.Range("K68").Value = 1396
$ws.Range("L68").Value = 1894.3334
$ws.Range("M68").Value = -647
$ws.Range("N68").Value = -3392.3334
$ws.Range("H71").Value = 1695
$ws.Range("I71").Value = 1396
$ws.Range("J71").Value = 1894.3334
$ws.Range("K71").Value = 6980
$ws.Range("L71").Value = 9471.666999999999
$ws.Range("M71").Value = -3236
$ws.Range("N71").Value = -16959.667
$ws.Range("H122").Value = 4226.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4226.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12679.9995
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17579.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 3533333.2
$ws.Range("I34").Value = 3533333.2
$ws.Range("K34").Value = 3533333.2
$ws.Range("M34").Value = -3533130.2
$ws.Range("H62").Value = 3000.9412
$ws.Range("I62").Value = 3169.111
$ws.Range("J62").Value = 2811.75
$ws.Range("K62").Value = 3169.111
$ws.Range("L62").Value = 2811.75
$ws.Range("M62").Value = -2545.111
$ws.Range("N62").Value = -4059.75
$ws.Range("H65").Value = 3000.9412
$ws.Range("I65").Value = 3169.111
$ws.Range("J65").Value = 2811.75
$ws.Range("K65").Value = 15845.555
$ws.Range("L65").Value = 14058.75
$ws.Range("M65").Value = -12725.555
$ws.Range("N65").Value = -20298.75
$ws.Range("H122").Value = 64998.125
$ws.Range("I122").Value = 113220
$ws.Range("K122").Value = 339660
$ws.Range("M122").Value = -337210
